# This script applies a permutation of data among rows 2,3,4,5,6,8,9 on the
# "Artfynd" sheet. The species-occurrence records (columns A,B,D,E,F,G,H,Q,R)
# are re-shuffled across those rows while column C ("Valideringsstatus") and
# all other columns stay untouched. Row 7 is unaffected.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values to write into each row, keyed by destination row number.
# Columns: A, B, D, E, F, G, H, Q, R
$rowsData = @{
    2 = @(97650299, 77506, "NT", 6425,   "Garnlav",             "Alectoria sarmentosa",  "(Ach.) Ach.",          403880.3826524244, 6794050.283030285)
    3 = @(97650294, 90676, "NT", 5966,   "Motaggsvamp",         "Sarcodon squamosus",    "(Schaeff.) Quél.",    403863.9880530759, 6794102.706117956)
    4 = @(97650292, 5135,  "LC", 105930, "Vågbandad barkbock",  "Semanotus undatus",     "(Linnaeus, 1758)",     403705.050704394,  6794737.908215457)
    5 = @(97650293, 77506, "NT", 6425,   "Garnlav",             "Alectoria sarmentosa",  "(Ach.) Ach.",          403710.6419448711, 6794695.894706693)
    6 = @(97650291, 5113,  "LC", 100526, "Bronshjon",           "Callidium coriaceum",   "Paykull, 1800",        403714.9324539425, 6794746.778207967)
    8 = @(97650301, 90653, "LC", 4364,   "Dropptaggsvamp",      "Hydnellum ferrugineum", "(Fr.:Fr.) P. Karst.", 403960.8920370748, 6793787.235077787)
    9 = @(97650298, 96334, "VU", 220787, "Knärot",              "Goodyera repens",       "(L.) R. Br.",          403840.5463236904, 6794038.864283022)
}

foreach ($r in $rowsData.Keys) {
    $vals = $rowsData[$r]

    $ws.Range("A$r").Value = $vals[0]
    $ws.Range("B$r").Value = $vals[1]
    $ws.Range("D$r").Value = $vals[2]
    $ws.Range("E$r").Value = $vals[3]
    $ws.Range("F$r").Value = $vals[4]
    $ws.Range("G$r").Value = $vals[5]
    $ws.Range("H$r").Value = $vals[6]
    $ws.Range("Q$r").Value = $vals[7]
    $ws.Range("R$r").Value = $vals[8]
}
